$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMP-103 Coding Adventures II")
$ws.Activate()

# Insert a new row at the top; existing rows shift down by one.
$ws.Range("A1").EntireRow.Insert()

# Fill in the new header row (write order matches shared-string insertion order).
$ws.Range("C1").Value2 = "Assignment 2"
$ws.Range("B1").Value2 = "Assignment 1"
$ws.Range("A1").Value2 = "Course Outcome"

# Bold the new header row.
$ws.Range("A1:C1").Font.Bold = $true

# Restore selection.
$ws.Range("A17").Select() | Out-Null
